# Update predicted-vs-actual ncv values (columns A, B, C) for rows 2-73
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 135
$ws.Range("B2").Value = 522
$ws.Range("C2").Value = 142
$ws.Range("A3").Value = 138
$ws.Range("B3").Value = 518
$ws.Range("C3").Value = 144
$ws.Range("A4").Value = 133
$ws.Range("B4").Value = 493
$ws.Range("C4").Value = 143
$ws.Range("A5").Value = 136
$ws.Range("B5").Value = 333
$ws.Range("C5").Value = 143
$ws.Range("A6").Value = 139
$ws.Range("B6").Value = 375
$ws.Range("C6").Value = 89
$ws.Range("A7").Value = 110
$ws.Range("B7").Value = 128
$ws.Range("C7").Value = 143
$ws.Range("A8").Value = 109
$ws.Range("B8").Value = 332
$ws.Range("C8").Value = 141
$ws.Range("A9").Value = 140
$ws.Range("B9").Value = 333
$ws.Range("C9").Value = 141
$ws.Range("A10").Value = 141
$ws.Range("B10").Value = 336
$ws.Range("C10").Value = 141
$ws.Range("A11").Value = 77
$ws.Range("B11").Value = 333
$ws.Range("C11").Value = 141
$ws.Range("A12").Value = 139
$ws.Range("B12").Value = 333
$ws.Range("C12").Value = 143
$ws.Range("A13").Value = 110
$ws.Range("B13").Value = 241
$ws.Range("C13").Value = 617
$ws.Range("A14").Value = 106
$ws.Range("B14").Value = 331
$ws.Range("C14").Value = 144
$ws.Range("A15").Value = 139
$ws.Range("B15").Value = 428
$ws.Range("C15").Value = 143
$ws.Range("A16").Value = 78
$ws.Range("B16").Value = 130
$ws.Range("C16").Value = 91
$ws.Range("A17").Value = 109
$ws.Range("B17").Value = 269
$ws.Range("C17").Value = 141
$ws.Range("A18").Value = 114
$ws.Range("B18").Value = 271
$ws.Range("C18").Value = 88
$ws.Range("A19").Value = 109
$ws.Range("B19").Value = 270
$ws.Range("C19").Value = 143
$ws.Range("A20").Value = 118
$ws.Range("B20").Value = 242
$ws.Range("C20").Value = 139
$ws.Range("A21").Value = 140
$ws.Range("B21").Value = 337
$ws.Range("C21").Value = 614
$ws.Range("A22").Value = 110
$ws.Range("B22").Value = 132
$ws.Range("C22").Value = 144
$ws.Range("A23").Value = 113
$ws.Range("B23").Value = 333
$ws.Range("C23").Value = 232
$ws.Range("A24").Value = 110
$ws.Range("B24").Value = 133
$ws.Range("C24").Value = 142
$ws.Range("A25").Value = 110
$ws.Range("B25").Value = 133
$ws.Range("C25").Value = 144
$ws.Range("A26").Value = 110
$ws.Range("B26").Value = 335
$ws.Range("C26").Value = 84
$ws.Range("A27").Value = 110
$ws.Range("B27").Value = 337
$ws.Range("C27").Value = 95
$ws.Range("A28").Value = 111
$ws.Range("B28").Value = 273
$ws.Range("C28").Value = 621
$ws.Range("A29").Value = 114
$ws.Range("B29").Value = 336
$ws.Range("C29").Value = 141
$ws.Range("A30").Value = 108
$ws.Range("B30").Value = 133
$ws.Range("C30").Value = 139
$ws.Range("A31").Value = 142
$ws.Range("B31").Value = 337
$ws.Range("C31").Value = 88
$ws.Range("A32").Value = 108
$ws.Range("B32").Value = 132
$ws.Range("C32").Value = 100
$ws.Range("A33").Value = 111
$ws.Range("B33").Value = 132
$ws.Range("C33").Value = 93
$ws.Range("A34").Value = 110
$ws.Range("B34").Value = 134
$ws.Range("C34").Value = 140
$ws.Range("A35").Value = 111
$ws.Range("B35").Value = 521
$ws.Range("C35").Value = 92
$ws.Range("A36").Value = 111
$ws.Range("B36").Value = 276
$ws.Range("C36").Value = 141
$ws.Range("A37").Value = 112
$ws.Range("B37").Value = 275
$ws.Range("C37").Value = 89
$ws.Range("A38").Value = 113
$ws.Range("B38").Value = 276
$ws.Range("C38").Value = 87
$ws.Range("A39").Value = 113
$ws.Range("B39").Value = 273
$ws.Range("C39").Value = 87
$ws.Range("A40").Value = 139
$ws.Range("B40").Value = 278
$ws.Range("C40").Value = 82
$ws.Range("A41").Value = 111
$ws.Range("B41").Value = 277
$ws.Range("C41").Value = 141
$ws.Range("A42").Value = 139
$ws.Range("B42").Value = 276
$ws.Range("C42").Value = 319
$ws.Range("A43").Value = 114
$ws.Range("B43").Value = 337
$ws.Range("C43").Value = 81
$ws.Range("A44").Value = 111
$ws.Range("B44").Value = 890
$ws.Range("C44").Value = 85
$ws.Range("A45").Value = 119
$ws.Range("B45").Value = 846
$ws.Range("C45").Value = 138
$ws.Range("A46").Value = 139
$ws.Range("B46").ClearContents() | Out-Null
$ws.Range("C46").Value = 140
$ws.Range("A47").Value = 138
$ws.Range("B47").ClearContents() | Out-Null
$ws.Range("C47").Value = 89
$ws.Range("A48").Value = 136
$ws.Range("B48").ClearContents() | Out-Null
$ws.Range("C48").Value = 86
$ws.Range("A49").Value = 138
$ws.Range("B49").ClearContents() | Out-Null
$ws.Range("C49").Value = 137
$ws.Range("A50").Value = 139
$ws.Range("B50").ClearContents() | Out-Null
$ws.Range("C50").Value = 137
$ws.Range("A51").Value = 116
$ws.Range("B51").ClearContents() | Out-Null
$ws.Range("C51").Value = 138
$ws.Range("A52").Value = 139
$ws.Range("C52").Value = 88
$ws.Range("A53").Value = 125
$ws.Range("C53").Value = 81
$ws.Range("A54").Value = 140
$ws.Range("C54").Value = 83
$ws.Range("A55").Value = 140
$ws.Range("C55").Value = 138
$ws.Range("A56").Value = 133
$ws.Range("C56").Value = 84
$ws.Range("A57").Value = 115
$ws.Range("C57").Value = 140
$ws.Range("A58").Value = 139
$ws.Range("C58").Value = 140
$ws.Range("A59").Value = 124
$ws.Range("C59").Value = 90
$ws.Range("A60").Value = 1011
$ws.Range("C60").Value = 91
$ws.Range("A61").Value = 137
$ws.Range("A62").Value = 130
$ws.Range("A63").Value = 138
$ws.Range("A64").Value = 139
$ws.Range("A65").Value = 144
$ws.Range("A66").Value = 141
$ws.Range("A67").Value = 138
$ws.Range("A68").Value = 137
$ws.Range("A69").Value = 140
$ws.Range("A70").Value = 138
$ws.Range("A71").Value = 138
$ws.Range("A72").Value = 140
$ws.Range("A73").Value = 148

# Update the visible selection to match the new data extent
$ws.Range("A1:C301").Select() | Out-Null
